# chore: update Sheets via scheduled runner
# Refreshes currentAveragePrice*/Leve profit figures (cols H-N) across the
# per-job leve sheets with newly pulled market data.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 831.63635
$ws.Range("J129").Value = 1000.55884
$ws.Range("L129").Value = 3001.67652
$ws.Range("N129").Value = -13001.67652
$ws.Range("H135").Value = 3887.9512
$ws.Range("I135").Value = 952.64
$ws.Range("J135").Value = 8474.375
$ws.Range("K135").Value = 8573.76
$ws.Range("L135").Value = 76269.375
$ws.Range("M135").Value = -6038.76
$ws.Range("N135").Value = -81339.375
$ws.Range("H137").Value = 4357.707
$ws.Range("I137").Value = 6881.346
$ws.Range("J137").Value = 2307.25
$ws.Range("K137").Value = 20644.038
$ws.Range("L137").Value = 6921.75
$ws.Range("M137").Value = -18094.038
$ws.Range("N137").Value = -12021.75
$ws.Range("H138").Value = 10641999
$ws.Range("I138").Value = 1667.1765
$ws.Range("J138").Value = 16671521
$ws.Range("K138").Value = 5001.529500000001
$ws.Range("L138").Value = 50014563
$ws.Range("M138").Value = 138.4704999999994
$ws.Range("N138").Value = -50024843

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null
$ws.Range("H36").Value = 51000
$ws.Range("I36").Value = 18666.666
$ws.Range("J36").Value = 83333.336
$ws.Range("K36").Value = 18666.666
$ws.Range("L36").Value = 83333.336
$ws.Range("M36").Value = -18320.666
$ws.Range("N36").Value = -84025.336
$ws.Range("H43").Value = 10792.333
$ws.Range("J43").Value = 10792.333
$ws.Range("L43").Value = 10792.333
$ws.Range("N43").Value = -11418.333
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2638
$ws.Range("H110").Value = 968.1
$ws.Range("J110").Value = 490
$ws.Range("L110").Value = 490
$ws.Range("N110").Value = -4580
$ws.Range("H132").Value = 166557.77
$ws.Range("I132").Value = 4024.92
$ws.Range("J132").Value = 843777.9399999999
$ws.Range("K132").Value = 12074.76
$ws.Range("L132").Value = 2531333.82
$ws.Range("M132").Value = -9544.76
$ws.Range("N132").Value = -2536393.82

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1365.5294
$ws.Range("I94").Value = 799.5
$ws.Range("K94").Value = 799.5
$ws.Range("M94").Value = -348.5
$ws.Range("H134").Value = 40774.78
$ws.Range("I134").Value = 49457.562
$ws.Range("J134").Value = 1219.8889
$ws.Range("K134").Value = 148372.686
$ws.Range("L134").Value = 3659.6667
$ws.Range("M134").Value = -145837.686
$ws.Range("N134").Value = -8729.6667

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 72249.75
$ws.Range("J9").Value = 72249.75
$ws.Range("L9").Value = 72249.75
$ws.Range("N9").Value = -72585.75
$ws.Range("H31").Value = 18641.582
$ws.Range("I31").Value = 27859.021
$ws.Range("J31").Value = 1624.7693
$ws.Range("K31").Value = 27859.021
$ws.Range("L31").Value = 1624.7693
$ws.Range("M31").Value = -27564.021
$ws.Range("N31").Value = -2214.7693
$ws.Range("H34").Value = 18641.582
$ws.Range("I34").Value = 27859.021
$ws.Range("J34").Value = 1624.7693
$ws.Range("K34").Value = 27859.021
$ws.Range("L34").Value = 1624.7693
$ws.Range("M34").Value = -27657.021
$ws.Range("N34").Value = -2028.7693
$ws.Range("H58").Value = 49049.285
$ws.Range("I58").Value = 1647.7693
$ws.Range("J58").Value = 126076.75
$ws.Range("K58").Value = 1647.7693
$ws.Range("L58").Value = 126076.75
$ws.Range("M58").Value = -1444.7693
$ws.Range("N58").Value = -126482.75
$ws.Range("H136").Value = 49049.285
$ws.Range("I136").Value = 1647.7693
$ws.Range("J136").Value = 126076.75
$ws.Range("K136").Value = 4943.3079
$ws.Range("L136").Value = 378230.25
$ws.Range("M136").Value = -2393.3079
$ws.Range("N136").Value = -383330.25

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 44199.715
$ws.Range("J37").Value = 44199.715
$ws.Range("L37").Value = 132599.145
$ws.Range("N37").Value = -132823.145
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -16372
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -51864
$ws.Range("H68").Value = 2244.9758
$ws.Range("I68").Value = 642.38635
$ws.Range("J68").Value = 4053.0256
$ws.Range("K68").Value = 1927.15905
$ws.Range("L68").Value = 12159.0768
$ws.Range("M68").Value = -1116.15905
$ws.Range("N68").Value = -13781.0768
$ws.Range("H71").Value = 2244.9758
$ws.Range("I71").Value = 642.38635
$ws.Range("J71").Value = 4053.0256
$ws.Range("K71").Value = 5781.47715
$ws.Range("L71").Value = 36477.2304
$ws.Range("M71").Value = -1725.47715
$ws.Range("N71").Value = -44589.2304
$ws.Range("H132").Value = 7640.684
$ws.Range("I132").Value = 6151
$ws.Range("J132").Value = 7815.9414
$ws.Range("K132").Value = 55359
$ws.Range("L132").Value = 70343.47259999999
$ws.Range("M132").Value = -52829
$ws.Range("N132").Value = -75403.47259999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 37292.75
$ws.Range("J62").Value = 37292.75
$ws.Range("L62").Value = 37292.75
$ws.Range("N62").Value = -38664.75
$ws.Range("H65").Value = 37292.75
$ws.Range("J65").Value = 37292.75
$ws.Range("L65").Value = 111878.25
$ws.Range("N65").Value = -118742.25
$ws.Range("H70").Value = 4397.278
$ws.Range("I70").Value = 4397.5
$ws.Range("K70").Value = 4397.5
$ws.Range("M70").Value = -4127.5
$ws.Range("H73").Value = 4397.278
$ws.Range("I73").Value = 4397.5
$ws.Range("K73").Value = 4397.5
$ws.Range("M73").Value = -3461.5
$ws.Range("H97").Value = 1271.1111
$ws.Range("J97").Value = 700
$ws.Range("L97").Value = 700
$ws.Range("N97").Value = -1692

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3761.3076
$ws.Range("I40").Value = 2389.2
$ws.Range("J40").Value = 8335
$ws.Range("K40").Value = 2389.2
$ws.Range("L40").Value = 8335
$ws.Range("M40").Value = -2253.2
$ws.Range("N40").Value = -8607
$ws.Range("H61").Value = 2121.2856
$ws.Range("I61").Value = 2609.8
$ws.Range("K61").Value = 2609.8
$ws.Range("M61").Value = -2407.8
$ws.Range("H63").Value = 46072.855
$ws.Range("J63").Value = 46072.855
$ws.Range("L63").Value = 46072.855
$ws.Range("N63").Value = -47570.855
$ws.Range("H66").Value = 46072.855
$ws.Range("J66").Value = 46072.855
$ws.Range("L66").Value = 138218.565
$ws.Range("N66").Value = -145706.565
$ws.Range("H113").Value = 2121.2856
$ws.Range("I113").Value = 2609.8
$ws.Range("K113").Value = 2609.8
$ws.Range("M113").Value = -439.8000000000002

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9740
$ws.Range("J74").Value = 9740
$ws.Range("L74").Value = 9740
$ws.Range("N74").Value = -11612
$ws.Range("H77").Value = 9740
$ws.Range("J77").Value = 9740
$ws.Range("L77").Value = 29220
$ws.Range("N77").Value = -38580
